# Revert "Bcrypt integration for secure password hashing - Part 2"
#
# The task list tracked "Use bcrypt for hashing and storing passwords*" (row 14)
# as already "Done", and had no entry for researching a password-hashing
# reference. Reverting that integration means:
#   - the bcrypt task's Status goes back to being un-set (blank)
#   - a new research-link task for crackstation's hashing article appears
#     at the bottom of the list
#   - two other in-flight tasks ("Add "Profile" page" and "Transition backend
#     to use database instead of static memory") are marked Done
#   - the view had scrolled down the list, leaving a new cell selected

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Add "Profile" page" (row 5): Status "In progress" -> "Done"
$ws.Range("D5").Value = "Done"

# "Use bcrypt for hashing and storing passwords*" (row 14): Status "Done" -> (blank)
$ws.Range("D14").Value = ""

# "Transition backend to use database instead of static memory" (row 15):
# Status "In progress" -> "Done"
$ws.Range("D15").Value = "Done"

# New task appended at row 31 (row 30 intentionally left empty), with only
# the Task column populated
$ws.Range("A31").Value = "https://crackstation.net/hashing-security.htm"

# Match the resulting selection left behind in the sheet view
$ws.Range("D21").Select() | Out-Null
